# "Append: 2025-09-07 12:30 JST" — the scraper re-ran and prepended two new
# listings to the "ランサーズ" sheet while re-stamping every existing row
# with the new retrieval timestamp.
#
# Resulting row layout (1-indexed, row 1 = header):
#   2  NEW   初回 AIヘルスケア...MVP開発パートナー募集
#   3  was 2 【募集】ジャーナリングとAIをテーマにしたiOSアプリ開発
#   4  NEW   【急募】年間カレンダー自動作成ツールの開発依頼
#   5  was 3 【急募】Instagram投稿を自動でGoogleビジネスに連携するMEOツール
#   6  was 4 IB報酬を得るための高性能EA開発依頼
#   7  was 5 限定公開 PR 限定公開の仕事

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the existing hyperlink entries before shuffling rows around — the
# underlying row-insert does not keep the worksheet's Hyperlinks collection
# in sync with the shifted cells, so we rebuild them from scratch afterwards.
$ws.Hyperlinks.Delete()

# Two brand-new rows land at 2 and 4. Insert them one at a time (the second
# insert happens *after* the first has already pushed everything down), which
# reproduces the observed before/after row mapping exactly.
$ws.Rows("2:2").Insert()
$ws.Rows("4:4").Insert()

$ts = "2025-09-07 12:30:54"

# Row 2 — new listing
$ws.Range("A2").Value = $ts
$ws.Range("B2").Value = "初回 AIヘルスケア×経営支援サービス|GPT-4・LINE API活用|MVP開発パートナー募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5388718"
$ws.Range("G2").Value = 635
$ws.Range("H2").Value = "🔥AI,GPT ◆開発"

# Row 3 — formerly row 2, timestamp refreshed only
$ws.Range("A3").Value = $ts

# Row 4 — new listing
$ws.Range("A4").Value = $ts
$ws.Range("B4").Value = "【急募】年間カレンダー自動作成ツールの開発依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5388837"
$ws.Range("G4").Value = 120
$ws.Range("H4").Value = "◆ツール,開発"

# Rows 5-7 — formerly rows 3-5, timestamp refreshed only
$ws.Range("A5").Value = $ts
$ws.Range("A6").Value = $ts
$ws.Range("A7").Value = $ts

# Rebuild hyperlinks on column F for every data row, in order, so relationship
# ids line up with the ref order the diff shows (rId1..rId6 == F2..F7).
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5388718")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5388502")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5388837")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5388589")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5388547")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5385681")

# Column B widened 42 -> 51 characters. The engine's save path adds a fixed
# +5/6 padding on top of whatever ColumnWidth is assigned, so back that out
# here to land on an exact integer width of 51 in the saved file.
$ws.Columns.Item(2).ColumnWidth = 51 - 5/6
